$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D text-forcing setup (avoid Excel auto-numeric conversion) ---
$dCells = @("D2", "D3", "D5", "D6", "D7", "D11", "D12", "D14", "D15", "D18", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($c in $dCells) { $ws.Range($c).NumberFormat = "@" }

# --- Apply new values ---
$ws.Range("D2").Value = "66.121.46"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "3.313.23"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "585.94"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "182.76"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +7.43%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").Value = "0.400"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "3.896.96"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "66.176.31"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "26.16"
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("D18").Value = "427.93"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D21").Value = "7.40"
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "71.91"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "5.69"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "3.462.45"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "0.516"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "0.204"
$ws.Range("E27").Value = "  +7.04%  "
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "8.92"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "1.95"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "22.28"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "5.17"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "6.58"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("D37").Value = "159.64"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").Value = "2.888.44"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "26.57"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "0.765"
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "4.31"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "40.13"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "0.0665"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("D46").Value = "5.97"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D50").Value = "0.0271"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "0.104"
$ws.Range("E51").Value = "  +4.39%  "

# --- Reset style on D cells so no stray style index is left on the cell ---
foreach ($c in $dCells) { $ws.Range($c).Style = "Normal" }

# --- Row swaps: rows that fully exchange coin identity + values ---
$swapDCells = @("D16", "D17", "D19", "D20", "D48", "D49")
foreach ($c in $swapDCells) { $ws.Range($c).NumberFormat = "@" }

# Row 16/17 swap: ShibaInu <-> WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.350.72"
$ws.Range("E16").Value = "  +0.32%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("E17").Value = "  -2.89%  "

# Row 19/20 swap: Chainlink <-> Polkadot
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "5.52"
$ws.Range("E19").Value = "  -3.01%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "13.18"
$ws.Range("E20").Value = "  -3.82%  "

# Row 48/49 swap: InjectiveProtocol <-> Bittensor
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "314.96"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "23.14"
$ws.Range("E49").Value = "  -6.62%  "

foreach ($c in $swapDCells) { $ws.Range($c).Style = "Normal" }
